$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had a stray "headless" row (row 13): it carried the
# professor's name in columns B/C but had no label in column A. That row is
# removed outright (whole-row delete), which shifts every row below it up by
# one — this is what actually happened in the authored edit (22 rows -> 21
# rows, uniqueCount 40 -> 37, count 59 -> 57).
$ws.Rows(13).Delete()

# From here on, row numbers refer to the NEW (post-delete) layout.

# "Objetivos:" (row 10) now shows the professor's identification instead of
# the long objectives paragraph.
$ws.Range("B10").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C10").Value = "9146830 - Danúbia Caporusso Bargos"

# "Programa resumido:" (row 13) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" (row 15) now repeats the same text used for "Ativação:" (row 8).
# Copy it cell-to-cell instead of retyping the literal text so Excel keeps
# storing it as shared text (matching the original style/shared-string
# layout) instead of auto-converting the "01/01/2018" string into a real
# date value/format.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# "Método:" (row 18) now shows the professor's identification too.
$ws.Range("B18").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C18").Value = "9146830 - Danúbia Caporusso Bargos"

# "Critério:" (row 19) now shows what used to be the "Método:" description.
$ws.Range("B19").Value = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."

# "Norma de recuperação:" (row 20) now shows what used to be the "Critério:"
# grading text.
$ws.Range("B20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."

# "Bibliografia:" (row 21) now shows what used to be the "Norma de
# recuperação:" text, and the old lengthy bibliography paragraph is dropped.
$ws.Range("B21").Value = "Provas e/ou exercícios dirigidos."
$ws.Range("C21").Value = "Provas e/ou exercícios dirigidos."
